$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Walmart"
$ws.Range("B13").Value = "AJH098987"
$ws.Range("C13").Value = "things"
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = "x"

$ws.Range("A14").Value = "Target"
$ws.Range("B14").Value = "OPJD000061"
$ws.Range("C14").Value = "other"
$ws.Range("D14").Value = 5
$ws.Range("F14").Value = "x"

$ws.Range("C16").Select() | Out-Null
